# Slide 52, "Content Placeholder 2": split the first paragraph's single
# run into three runs so that the sentence reads "...what we did for
# global variables..." (inserting the word "for ").
#
#   Before: "Similar to what we did global variables, ..."
#   After:  "Similar to what " + "we did for " + "global variables, ..."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(52)
$shape = $s.Shapes.Item(2)

$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(1, 1)

# Characters 17-23 of the paragraph are "we did " (1-based, inclusive).
# Replacing that span with "we did for " splits the original single run
# into the three runs required by the edit, with the new word "for "
# folded into the (previously existing) middle segment.
$mid = $para.Characters(17, 7)
$mid.Text = "we did for "
